# Updates "cryptos list" (Price / Volume(1h)) values for rows 2-51 on Sheet1.
# D = Price (col 4), E = Volume(1h) (col 5). Values that read as plain numbers
# are pre-formatted as Text ("@") before assignment so Excel keeps them as text
# (matching the source workbook, where these cells are inline/shared strings,
# e.g. "29.002.84" or "1.014", not numeric values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "28.989.28"; E = "  -1.73%  " }
    @{ Row = 3; D = "1.968.00"; E = "  -1.33%  " }
    @{ Row = 4; D = $null; E = "  +0.78%  " }
    @{ Row = 5; D = "328.91"; E = "  -0.19%  " }
    @{ Row = 6; D = "1.014"; E = "  +0.68%  " }
    @{ Row = 7; D = "0.4952"; E = "  -1.07%  " }
    @{ Row = 8; D = "0.4171"; E = "  -1.10%  " }
    @{ Row = 9; D = "54.42"; E = "  +4.65%  " }
    @{ Row = 10; D = "0.09214"; E = "  +3.47%  " }
    @{ Row = 11; D = "1.089"; E = "  -2.82%  " }
    @{ Row = 12; D = "22.61"; E = "  -3.14%  " }
    @{ Row = 13; D = "1.999.07"; E = "  +0.75%  " }
    @{ Row = 14; D = "7.844"; E = "  -3.10%  " }
    @{ Row = 15; D = "6.425"; E = "  -1.42%  " }
    @{ Row = 16; D = "1.017"; E = "  +0.91%  " }
    @{ Row = 17; D = "0.00001106"; E = "  -0.23%  " }
    @{ Row = 18; D = "91.33"; E = "  -5.06%  " }
    @{ Row = 19; D = "0.06733"; E = "  +1.47%  " }
    @{ Row = 20; D = "19.05"; E = "  -3.46%  " }
    @{ Row = 21; D = $null; E = "  +0.67%  " }
    @{ Row = 22; D = "5.949"; E = "  -0.31%  " }
    @{ Row = 23; D = "28.999.41"; E = "  -1.77%  " }
    @{ Row = 24; D = "11.90"; E = $null }
    @{ Row = 25; D = "2.273"; E = "  +0.04%  " }
    @{ Row = 26; D = "2.259.96"; E = "  -0.21%  " }
    @{ Row = 27; D = "20.67"; E = "  +0.23%  " }
    @{ Row = 28; D = "156.33"; E = "  -1.13%  " }
    @{ Row = 29; D = "6.201"; E = "  -5.83%  " }
    @{ Row = 30; D = "2.250"; E = "  -3.65%  " }
    @{ Row = 31; D = "126.86"; E = "  -0.83%  " }
    @{ Row = 32; D = "1.038"; E = "  -1.18%  " }
    @{ Row = 33; D = "0.09795"; E = "  -1.50%  " }
    @{ Row = 34; D = "1.498"; E = "  -3.87%  " }
    @{ Row = 35; D = "5.789"; E = "  -0.83%  " }
    @{ Row = 36; D = "3.753"; E = "  -0.91%  " }
    @{ Row = 37; D = "0.02407"; E = "  -2.24%  " }
    @{ Row = 38; D = "1.312"; E = "  +1.72%  " }
    @{ Row = 39; D = "0.06365"; E = "  +0.05%  " }
    @{ Row = 40; D = "8.991"; E = "  -6.08%  " }
    @{ Row = 41; D = "0.6437"; E = "  -1.40%  " }
    @{ Row = 42; D = "11.38"; E = "  -2.90%  " }
    @{ Row = 43; D = "0.1990"; E = "  -3.70%  " }
    @{ Row = 44; D = "1.014"; E = "  +0.75%  " }
    @{ Row = 45; D = "0.6154"; E = "  -3.02%  " }
    @{ Row = 46; D = "1.339"; E = "  +5.21%  " }
    @{ Row = 47; D = "13.20"; E = "  -1.29%  " }
    @{ Row = 48; D = "2.158"; E = "  -2.44%  " }
    @{ Row = 49; D = "3.487"; E = "  -1.07%  " }
    @{ Row = 50; D = "0.00000000332"; E = "  +0.91%  " }
    @{ Row = 51; D = "0.06941"; E = "  -0.95%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
